# Update the "A5" sheet's ClassID/ChainID data row, leaving the cursor
# parked on B15 (no longer the selected tab).
$wb = $excel.ActiveWorkbook

$ws5 = $wb.Worksheets.Item("A5")
$ws5.Range("B2").Value = "juno1wrq62rjy9w07lxz683kdsrrpqa2mnpkf2dd7hm4zhzh4pld76ems2w3rcx"
$ws5.Range("D2").Value = "uni6"
$ws5.Range("B15").Select()

# Update the "A6" sheet's ClassID/ChainID data row, then make it the
# active/selected tab with A2 selected.
$ws6 = $wb.Worksheets.Item("A6")
$ws6.Range("B2").Value = "ibc/FFD8922D35939D31F1196AFF94211EDC5AD9E36117CF5464FA4B56530B4236E5"
$ws6.Range("D2").Value = "gon-flixnet-1"
$ws6.Activate()
$ws6.Range("A2").Select()
